$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date value (serial 45204 = 2023-10-05) for every
# data row (rows 2..289). The update bumps that date by one day (45205 = 2023-10-06)
# for all of them.
$ws.Range("C2:C289").Value = 45205
